$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '320.39'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '7.60%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '48.93'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '17.13%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.265'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '4.96%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08104'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '7.65%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.609'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '5.37%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.663'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.58%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.200'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '30.72%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '11.67%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1945'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '6.53%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09518'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '6.23%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04510'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '9.70%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1049'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.12%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001328'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '3.38%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005947'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.56%'
$ws.Range("B16").Value = 'HotbitToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.004246'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '8.12%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.358'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.58%'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.436'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.46%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3391'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.89%'
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.197'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.06%'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1412'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '3.00%'
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.3061'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-4.95%'
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04298'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '5.01%'
$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001311'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '3.54%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001353'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '3.98%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003547'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-4.74%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02680'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '11.35%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05572'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '7.18%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.006314'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.16%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007700'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.41%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1438'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '8.56%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007711'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4.21%'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '14.09%'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.84%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007015'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '6.47%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000752'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.16%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06083'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '30.04%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004009'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-4.61%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002105'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.16%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002004'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.16%'
